$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) from "Gamma1F-HW50.xpc" to "Gamma1F"
$ws.Name = "Gamma1F"

# Add a new row of data (row 16) following the same pattern as the
# previous rows: A = index (14), B = category label (same text as row 15),
# C:M = 1
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1

# Match the formatting style of column A used on previous rows (index column)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
